# New crime data collected - weekly CompStat update (105th Precinct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# "Volume 31   Number  10" -> "... 11"
$ws.Range("A8").Value = "Volume 31   Number  11"
# "Report Covering the Week  3/4/2024  Through  3/10/2024" -> new week
$ws.Range("C9").Value = "Report Covering the Week  3/11/2024  Through  3/17/2024"

# --- Row 15 (Rape) ---------------------------------------------------------
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 0

# --- Row 16 (Robbery) -------------------------------------------------------
# C16 switches from a numeric cell to the text placeholder "0" (same as C14/C22)
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -57.142857142857
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = 7.407407407407
$ws.Range("L16").Value = -12.121212121212
$ws.Range("M16").Value = -63.75
$ws.Range("N16").Value = -86.320754716981

# --- Row 17 (Fel. Assault) --------------------------------------------------
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -3.333333333333
$ws.Range("I17").Value = 98
$ws.Range("J17").Value = 85
$ws.Range("K17").Value = 15.294117647058
$ws.Range("L17").Value = 10.112359550561
$ws.Range("M17").Value = 75
$ws.Range("N17").Value = 40

# --- Row 18 (Burglary) ------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = -51.724137931034
$ws.Range("I18").Value = 35
$ws.Range("J18").Value = 48
$ws.Range("K18").Value = -27.083333333333
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -53.333333333333
$ws.Range("N18").Value = -87.132352941176

# --- Row 19 (Gr. Larceny) ---------------------------------------------------
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 57.142857142857
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = -2.380952380952
$ws.Range("I19").Value = 131
$ws.Range("J19").Value = 121
$ws.Range("K19").Value = 8.264462809917
$ws.Range("L19").Value = 12.931034482758
$ws.Range("M19").Value = 48.863636363636
$ws.Range("N19").Value = 29.702970297029

# --- Row 20 (G.L.A.) ---------------------------------------------------------
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 17
$ws.Range("H20").Value = -37.037037037037
$ws.Range("I20").Value = 70
$ws.Range("J20").Value = 58
$ws.Range("K20").Value = 20.689655172413
$ws.Range("L20").Value = 45.833333333333
$ws.Range("M20").Value = -10.256410256410
$ws.Range("N20").Value = -89.567809239940

# --- Row 21 (TOTAL) ----------------------------------------------------------
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 14.814814814814
$ws.Range("G21").Value = 136
$ws.Range("H21").Value = -21.323529411764
$ws.Range("I21").Value = 370
$ws.Range("J21").Value = 342
$ws.Range("K21").Value = 8.187134502923
$ws.Range("L21").Value = 13.846153846153
$ws.Range("M21").Value = -4.145077720207
$ws.Range("N21").Value = -72.284644194756

# --- Row 24 (Petit Larceny) ---------------------------------------------------
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 5.882352941176
$ws.Range("F24").Value = 79
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = -24.761904761904
$ws.Range("I24").Value = 246
$ws.Range("J24").Value = 275
$ws.Range("K24").Value = -10.545454545454
$ws.Range("L24").Value = -12.142857142857
$ws.Range("M24").Value = 57.692307692307

# --- Row 25 (Retail Theft) ----------------------------------------------------
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 50
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 36.363636363636
$ws.Range("I25").Value = 48
$ws.Range("J25").Value = 46
$ws.Range("K25").Value = 4.347826086956
$ws.Range("L25").Value = -11.111111111111

# --- Row 26 (Misd. Assault) ---------------------------------------------------
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 64
$ws.Range("G26").Value = 41
$ws.Range("H26").Value = 56.097560975609
$ws.Range("I26").Value = 151
$ws.Range("J26").Value = 113
$ws.Range("K26").Value = 33.628318584070
$ws.Range("L26").Value = 32.456140350877
$ws.Range("M26").Value = 9.420289855072

# --- Row 27 (UCR Rape*) --------------------------------------------------------
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 9
$ws.Range("K27").Value = 28.571428571428
$ws.Range("L27").Value = -10

# --- Row 28 (Other Sex Crimes) -------------------------------------------------
# D28 switches from the text placeholder "0" to a real number
$ws.Range("D28").Value = 1
$ws.Range("D16").Copy()
$ws.Range("D28").PasteSpecial(-4122)
# E28 switches from the text placeholder "***.*" to a real percentage number
$ws.Range("E28").Value = 0
$ws.Range("E16").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 11
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = 37.5
$ws.Range("L28").Value = 83.333333333333

# --- Row 29 (Shooting Vic.) -----------------------------------------------------
# C29 switches from the text placeholder "0" to a real number
$ws.Range("C29").Value = 1
$ws.Range("F29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("F29").Value = 2
$ws.Range("I29").Value = 3
$ws.Range("K29").Value = 50
$ws.Range("L29").Value = 50
$ws.Range("N29").Value = -66.666666666666

# --- Row 30 (Shooting Inc.) ------------------------------------------------------
# C30 switches from the text placeholder "0" to a real number
$ws.Range("C30").Value = 1
$ws.Range("F30").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 3
$ws.Range("K30").Value = 200
$ws.Range("L30").Value = 200
$ws.Range("M30").Value = -57.142857142857
$ws.Range("N30").Value = -66.666666666666
